$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $text) {
    # Column D ("Price") holds numeric-looking strings (e.g. "602.22",
    # "68.560.73", "0.0000293") that must stay literal text, exactly as
    # authored, instead of being auto-coerced into floating point numbers
    # by the usual COM "smart" type detection. A leading apostrophe forces
    # Excel to store the literal text (same as a user typing '602.22 into
    # the cell); resetting the Style back to "Normal" afterwards keeps the
    # cell's formatting identical to its original (unset) state.
    $cell = $ws.Range($addr)
    $cell.Value = "'" + $text
    $cell.Style = "Normal"
}

# Row 2 - Bitcoin
Set-TextValue "D2" "68.560.73"
$ws.Range("E2").Value = "  -2.66%  "

# Row 3 - Ethereum
Set-TextValue "D3" "3.709.86"
$ws.Range("E3").Value = "  -3.32%  "

# Row 4 - TetherUSD
Set-TextValue "D4" "0.999"
$ws.Range("E4").Value = "  -0.24%  "

# Row 5 - BNB
Set-TextValue "D5" "602.22"
$ws.Range("E5").Value = "  +1.79%  "

# Row 6 - Solana
Set-TextValue "D6" "183.49"
$ws.Range("E6").Value = "  +9.86%  "

# Row 7 - LidoStakedEther
Set-TextValue "D7" "3.700.23"
$ws.Range("E7").Value = "  -3.43%  "

# Row 8 - XRP
$ws.Range("E8").Value = "  -5.62%  "

# Row 9 - USDC
Set-TextValue "D9" "0.997"
$ws.Range("E9").Value = "  -0.32%  "

# Row 10 - Cardano
Set-TextValue "D10" "0.719"
$ws.Range("E10").Value = "  -3.38%  "

# Row 11 - Dogecoin
$ws.Range("E11").Value = "  -6.22%  "

# Row 12 - Avalanche
Set-TextValue "D12" "56.71"
$ws.Range("E12").Value = "  +7.28%  "

# Row 13 - ShibaInu
Set-TextValue "D13" "0.0000293"
$ws.Range("E13").Value = "  -8.22%  "

# Row 14 - Polkadot
Set-TextValue "D14" "10.45"
$ws.Range("E14").Value = "  -7.51%  "

# Row 15 - WrappedliquidstakedEther2.0
Set-TextValue "D15" "4.291.48"
$ws.Range("E15").Value = "  -3.81%  "

# Row 16 - WrappedEther
Set-TextValue "D16" "3.711.53"
$ws.Range("E16").Value = "  -3.52%  "

# Row 17 - Chainlink
Set-TextValue "D17" "19.40"
$ws.Range("E17").Value = "  -7.47%  "

# Row 18 - TRON
$ws.Range("E18").Value = "  -2.03%  "

# Row 19 - Uniswap
Set-TextValue "D19" "12.90"
$ws.Range("E19").Value = "  -6.29%  "

# Row 20 - Polygon
Set-TextValue "D20" "1.13"
$ws.Range("E20").Value = "  -6.38%  "

# Row 21 - WrappedBTC
Set-TextValue "D21" "68.198.30"
$ws.Range("E21").Value = "  -3.17%  "

# Row 22 - BitcoinCash
Set-TextValue "D22" "411.41"
$ws.Range("E22").Value = "  -5.42%  "

# Row 23 - PancakeSwap
$ws.Range("E23").Value = "  -1.22%  "

# Row 24 - Litecoin
Set-TextValue "D24" "89.37"
$ws.Range("E24").Value = "  -4.65%  "

# Row 25 - ImmutableX
$ws.Range("E25").Value = "  -6.51%  "

# Row 26 - InternetComputer(DFINITY)
Set-TextValue "D26" "12.89"
$ws.Range("E26").Value = "  -6.82%  "

# Row 27 - was Toncoin, now RenderToken
$ws.Range("B27").Value = "RenderToken"
$ws.Range("C27").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue "D27" "10.92"
$ws.Range("E27").Value = "  -1.94%  "

# Row 28 - was RenderToken, now Toncoin
$ws.Range("B28").Value = "Toncoin"
$ws.Range("C28").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
Set-TextValue "D28" "3.92"
$ws.Range("E28").Value = "  -2.12%  "

# Row 29 - LEO
$ws.Range("E29").Value = "  +1.83%  "

# Row 30 - Filecoin
Set-TextValue "D30" "9.50"
$ws.Range("E30").Value = "  -8.49%  "

# Row 31 - EthereumClassic
Set-TextValue "D31" "32.98"
$ws.Range("E31").Value = "  -5.65%  "

# Row 32 - NEARProtocol
Set-TextValue "D32" "7.30"
$ws.Range("E32").Value = "  -9.99%  "

# Row 33 - Cosmos
Set-TextValue "D33" "12.52"
$ws.Range("E33").Value = "  -6.52%  "

# Row 34 - Hedera
Set-TextValue "D34" "0.118"
$ws.Range("E34").Value = "  -5.58%  "

# Row 35 - InjectiveProtocol
Set-TextValue "D35" "43.86"
$ws.Range("E35").Value = "  -8.16%  "

# Row 36 - OKB
Set-TextValue "D36" "64.74"
$ws.Range("E36").Value = "  -6.57%  "

# Row 37 - Bittensor
Set-TextValue "D37" "602.53"
$ws.Range("E37").Value = "  -4.96%  "

# Row 38 - PEPE (value contains a subscript-3 character, U+2083)
Set-TextValue "D38" "0.0₃0894"
$ws.Range("E38").Value = "  -8.48%  "

# Row 39 - was Dai, now TheGraph
$ws.Range("B39").Value = "TheGraph"
$ws.Range("C39").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
Set-TextValue "D39" "0.401"
$ws.Range("E39").Value = "  -4.93%  "

# Row 40 - was TheGraph, now Dai
$ws.Range("B40").Value = "Dai"
$ws.Range("C40").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
Set-TextValue "D40" "1.00"
$ws.Range("E40").Value = "  +0.04%  "

# Row 41 - FirstDigitalUSD
$ws.Range("E41").Value = "  +0.08%  "

# Row 42 - Kaspa
Set-TextValue "D42" "0.137"
$ws.Range("E42").Value = "  -5.72%  "

# Row 43 - ThetaToken
Set-TextValue "D43" "3.06"
$ws.Range("E43").Value = "  -6.13%  "

# Row 44 - Fetch.AI
Set-TextValue "D44" "2.76"
$ws.Range("E44").Value = "  +2.16%  "

# Row 45 - VeChain
Set-TextValue "D45" "0.0441"
$ws.Range("E45").Value = "  -5.65%  "

# Row 46 - dogwifhat
Set-TextValue "D46" "2.89"
$ws.Range("E46").Value = "  -10.05%  "

# Row 47 - THORChain
Set-TextValue "D47" "9.25"
$ws.Range("E47").Value = "  -7.47%  "

# Row 48 - WEMIXToken
Set-TextValue "D48" "2.75"
$ws.Range("E48").Value = "  -2.95%  "

# Row 49 - was Maker, now Stellar
$ws.Range("B49").Value = "Stellar"
$ws.Range("C49").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
Set-TextValue "D49" "0.135"
$ws.Range("E49").Value = "  -5.35%  "

# Row 50 - was Stellar, now Maker
$ws.Range("B50").Value = "Maker"
$ws.Range("C50").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
Set-TextValue "D50" "2.788.29"
$ws.Range("E50").Value = "  -2.05%  "

# Row 51 - ApeXProtocol
$ws.Range("E51").Value = "  -3.02%  "
